# Scheduled runner update: refresh computed Leve profit columns (H-N)
# for rows whose market-board prices changed since the last sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 134 (Leve Item ID / G134 = 41997)
$ws.Range("H134").Value = 69122.22
$ws.Range("J134").Value = 69122.22
$ws.Range("L134").Value = 69122.22
$ws.Range("N134").Value = -79262.22

# Row 138 (Leve Item ID / G138 = 44169)
$ws.Range("H138").Value = 129447.85
$ws.Range("I138").Value = 2115.0454
$ws.Range("J138").Value = 173218.5
$ws.Range("K138").Value = 6345.1362
$ws.Range("L138").Value = 519655.5
$ws.Range("M138").Value = -1205.1362
$ws.Range("N138").Value = -529935.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID / G45 = 27714)
$ws.Range("H45").Value = 3259.261
$ws.Range("I45").Value = 3294.3333
$ws.Range("J45").Value = 3236.7144
$ws.Range("K45").Value = 3294.3333
$ws.Range("L45").Value = 3236.7144
$ws.Range("M45").Value = -2917.3333
$ws.Range("N45").Value = -3990.7144

# Row 74 (Leve Item ID / G74 = 44000)
$ws.Range("H74").Value = 1258.96
$ws.Range("I74").Value = 974.9524
$ws.Range("J74").Value = 2750
$ws.Range("K74").Value = 974.9524
$ws.Range("L74").Value = 2750
$ws.Range("M74").Value = -100.9524
$ws.Range("N74").Value = -4498

# Row 77 (Leve Item ID / G77 = 44000)
$ws.Range("H77").Value = 1258.96
$ws.Range("I77").Value = 974.9524
$ws.Range("J77").Value = 2750
$ws.Range("K77").Value = 4874.762
$ws.Range("L77").Value = 13750
$ws.Range("M77").Value = -506.7619999999997
$ws.Range("N77").Value = -22486

# Row 122 (Leve Item ID / G122 = 36168)
$ws.Range("H122").Value = 1964.4615
$ws.Range("I122").Value = 2146.2856
$ws.Range("J122").Value = 1752.3334
$ws.Range("K122").Value = 6438.8568
$ws.Range("L122").Value = 5257.0002
$ws.Range("M122").Value = -3988.8568
$ws.Range("N122").Value = -10157.0002

# Row 123 (Leve Item ID / G123 = 34107)
$ws.Range("H123").Value = 33428
$ws.Range("J123").Value = 33428
$ws.Range("L123").Value = 33428
$ws.Range("N123").Value = -43228

# Row 132 (Leve Item ID / G132 = 43997)
$ws.Range("H132").Value = 4300.4517
$ws.Range("I132").Value = 3263.1667
$ws.Range("J132").Value = 7856.857
$ws.Range("K132").Value = 9789.500100000001
$ws.Range("L132").Value = 23570.571
$ws.Range("M132").Value = -7259.500100000001
$ws.Range("N132").Value = -28630.571

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (Leve Item ID / G107 = 27706)
$ws.Range("H107").Value = 1400.4286
$ws.Range("I107").Value = 1189.5555
$ws.Range("J107").Value = 1780
$ws.Range("K107").Value = 1189.5555
$ws.Range("L107").Value = 1780
$ws.Range("M107").Value = 730.4445000000001
$ws.Range("N107").Value = -5620

# Row 134 (Leve Item ID / G134 = 43998)
$ws.Range("H134").Value = 3235.8235
$ws.Range("I134").Value = 3012.2307
$ws.Range("J134").Value = 3962.5
$ws.Range("K134").Value = 9036.6921
$ws.Range("L134").Value = 11887.5
$ws.Range("M134").Value = -6501.6921
$ws.Range("N134").Value = -16957.5

$ws = $wb.Worksheets.Item("CRP")
# Row 2 (Leve Item ID / G2 = 1820)
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 22 (Leve Item ID / G22 = 5367)
$ws.Range("H22").Value = 328.27274
$ws.Range("I22").Value = 234.55556
$ws.Range("K22").Value = 234.55556
$ws.Range("M22").Value = 115.44444

# Row 31 (Leve Item ID / G31 = 44023)
$ws.Range("H31").Value = 8371.931
$ws.Range("I31").Value = 2034.1818
$ws.Range("J31").Value = 12245
$ws.Range("K31").Value = 2034.1818
$ws.Range("L31").Value = 12245
$ws.Range("M31").Value = -1739.1818
$ws.Range("N31").Value = -12835

# Row 34 (Leve Item ID / G34 = 44023)
$ws.Range("H34").Value = 8371.931
$ws.Range("I34").Value = 2034.1818
$ws.Range("J34").Value = 12245
$ws.Range("K34").Value = 2034.1818
$ws.Range("L34").Value = 12245
$ws.Range("M34").Value = -1832.1818
$ws.Range("N34").Value = -12649

# Row 105 (Leve Item ID / G105 = 19928)
$ws.Range("H105").Value = 600
$ws.Range("I105").Value = 550
$ws.Range("K105").Value = 550
$ws.Range("M105").Value = 1197

$ws = $wb.Worksheets.Item("CUL")
# Row 121 (Leve Item ID / G121 = 27878)
$ws.Range("H121").Value = 1037.1666
$ws.Range("I121").Value = 600.9091
$ws.Range("J121").Value = 1135.102
$ws.Range("K121").Value = 1802.7273
$ws.Range("L121").Value = 3405.306
$ws.Range("M121").Value = -492.7273
$ws.Range("N121").Value = -6025.306

$ws = $wb.Worksheets.Item("GSM")
# Row 15 (Leve Item ID / G15 = 12018)
$ws.Range("H15").Value = 34000
$ws.Range("J15").Value = 34000
$ws.Range("L15").Value = 34000
$ws.Range("N15").Value = -34576

# Row 63 (Leve Item ID / G63 = 11048)
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372

# Row 66 (Leve Item ID / G66 = 11048)
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864

# Row 81 (Leve Item ID / G81 = 12018)
$ws.Range("H81").Value = 34000
$ws.Range("J81").Value = 34000
$ws.Range("L81").Value = 34000
$ws.Range("N81").Value = -35996

# Row 84 (Leve Item ID / G84 = 12018)
$ws.Range("H84").Value = 34000
$ws.Range("J84").Value = 34000
$ws.Range("L84").Value = 102000
$ws.Range("N84").Value = -111984

# Row 118 (Leve Item ID / G118 = 26172)
$ws.Range("H118").Value = 42560
$ws.Range("J118").Value = 42560
$ws.Range("L118").Value = 42560
$ws.Range("N118").Value = -45874

# Row 121 (Leve Item ID / G121 = 26338)
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# Row 124 (Leve Item ID / G124 = 34247)
$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws = $wb.Worksheets.Item("LTW")
# Row 43 (Leve Item ID / G43 = 4314)
$ws.Range("H43").Value = 27678.4
$ws.Range("J43").Value = 27678.4
$ws.Range("L43").Value = 27678.4
$ws.Range("N43").Value = -28064.4

# Row 131 (Leve Item ID / G131 = 35466)
$ws.Range("H131").Value = 9000
$ws.Range("J131").Value = 9000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

# Row 132 (Leve Item ID / G132 = 44058)
$ws.Range("H132").Value = 4288.4585
$ws.Range("I132").Value = 3902.1333
$ws.Range("J132").Value = 4932.3335
$ws.Range("K132").Value = 11706.3999
$ws.Range("L132").Value = 14797.0005
$ws.Range("M132").Value = -9176.3999
$ws.Range("N132").Value = -19857.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (Leve Item ID / G122 = 36208)
$ws.Range("H122").Value = 1643.25
$ws.Range("I122").Value = 1692.8
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 5078.4
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -2628.4
$ws.Range("N122").Value = -7600

# Row 123 (Leve Item ID / G123 = 34127)
$ws.Range("H123").Value = 46275.75
$ws.Range("J123").Value = 46275.75
$ws.Range("L123").Value = 46275.75
$ws.Range("N123").Value = -56075.75

# Row 132 (Leve Item ID / G132 = 44029)
$ws.Range("H132").Value = 11908235
$ws.Range("I132").Value = 4335.3335
$ws.Range("J132").Value = 20836160
$ws.Range("K132").Value = 13006.0005
$ws.Range("L132").Value = 62508480
$ws.Range("M132").Value = -10476.0005
$ws.Range("N132").Value = -62513540
